$d = $word.ActiveDocument

function Set-ParagraphXml {
    param(
        [int]$ParaIndex,
        [string]$InnerXml
    )

    $p = $d.Paragraphs($ParaIndex)
    $r = $p.Range
    $pStart = $r.Start
    $textLen = $r.Text.Length

    # Clear all run text in the paragraph but keep the paragraph mark itself,
    # so the paragraph (and its identity/properties) survives.
    if ($textLen -gt 1) {
        $clearRange = $d.Range($pStart, $pStart + $textLen - 1)
        $clearRange.Text = ""
    }

    $p2 = $d.Paragraphs($ParaIndex)
    $insertRange = $d.Range($p2.Range.Start, $p2.Range.Start)

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $InnerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insertRange.InsertXML($xml)
}

# Paragraph 5: "We finally reach the end and exit the loop.  We then return the array value stored at index -m. "
# becomes: "We finally reach the end and exit the loop, storing the length. We then iterate over the
#           nodes again, and this time stopping at our target node which is length – m + 1.  We return that node."
$para5Xml = '<w:r><w:t>W</w:t></w:r>' +
            '<w:r><w:t>e finally reach the end</w:t></w:r>' +
            '<w:r><w:t xml:space="preserve"> and </w:t></w:r>' +
            '<w:r><w:t>exit the loop</w:t></w:r>' +
            '<w:r><w:t>, storing the length. We then iterate over the nodes again, and this time stopping at our target node which is length – m + 1.  We return that node.</w:t></w:r>'
Set-ParagraphXml 5 $para5Xml

# Paragraph 6: "Time efficiency of this is O(N) since we only need to traverse the list once. Space efficiency is O(N) as well."
# becomes: "Time efficiency of this is O(2n) since we traverse the list twice. Space efficiency is negligible
#           since we are only using small variables."
$para6Xml = '<w:r><w:t>Time e</w:t></w:r>' +
            '<w:r><w:t>fficiency of this is O(</w:t></w:r>' +
            '<w:r><w:t>2</w:t></w:r>' +
            '<w:r><w:t>n</w:t></w:r>' +
            '<w:r><w:t xml:space="preserve">) since we traverse the list </w:t></w:r>' +
            '<w:r><w:t>twice</w:t></w:r>' +
            '<w:r><w:t>.</w:t></w:r>' +
            '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
            '<w:r><w:t>Space efficiency is negligible since we are only using small variables.</w:t></w:r>' +
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
Set-ParagraphXml 6 $para6Xml

Write-Host "Done"
